# "small changes to slides"
#
# Slide 6 ("First WEek", sldId 260) has its "Content Placeholder 2" shape
# (id=3) text body tweaked:
#   - "Installing Software"  -> "Installing Software:"
#   - "And on windows"       -> "And on windows (you will not need these if you are on a mac):"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$hit1 = $tr.Find("Installing Software")
$hit1.Text = "Installing Software:"

$hit2 = $tr.Find("And on windows")
$hit2.Text = "And on windows (you will not need these if you are on a mac):"
